$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 307, pushing the existing rows 307:383 down to 308:384
$ws.Rows.Item(307).Insert()

# Fill the new row 307 with the new record's data (same constant columns as the rest of this block)
$ws.Cells.Item(307, 1).Value = 7
$ws.Cells.Item(307, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(307, 3).Value = "Ñuble"
$ws.Cells.Item(307, 4).Value = 44943
$ws.Cells.Item(307, 5).Value = 16
$ws.Cells.Item(307, 6).Value = 100112002
$ws.Cells.Item(307, 7).Value = "Pimiento"
$ws.Cells.Item(307, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(307, 9).Value = "Primera"
$ws.Cells.Item(307, 10).Value = 60
$ws.Cells.Item(307, 11).Value = 10000
$ws.Cells.Item(307, 12).Value = 10000
$ws.Cells.Item(307, 13).Value = 10000
$ws.Cells.Item(307, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(307, 15).Value = "Región del Maule"
$ws.Cells.Item(307, 16).Value = 556
$ws.Cells.Item(307, 17).Value = 18
$ws.Cells.Item(307, 18).Value = "Hortaliza"
